$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter header labels (row 8 and row 24) ---
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل اول منتهی به 1401/12"

$ws.Range("E24").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F24").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G24").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H24").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I24").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J24").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K24").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L24").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M24").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N24").Value = "فصل اول منتهی به 1401/12"

# --- Update quarterly data values ---
$ws.Range("E10").Value = 26992
$ws.Range("F10").Value = 38871
$ws.Range("G10").Value = 20143
$ws.Range("H10").Value = 60972
$ws.Range("I10").Value = 59010
$ws.Range("J10").Value = 35106
$ws.Range("K10").Value = 34932
$ws.Range("L10").Value = 44491
$ws.Range("M10").Value = 140564
$ws.Range("N10").Value = 65892

$ws.Range("E13").Value = 115210
$ws.Range("F13").Value = 67562
$ws.Range("G13").Value = 57286
$ws.Range("H13").Value = 292585
$ws.Range("I13").Value = 50733
$ws.Range("J13").Value = 47468
$ws.Range("K13").Value = -8619
$ws.Range("L13").Value = 20265
$ws.Range("M13").Value = 13106
$ws.Range("N13").Value = 11315

$ws.Range("E14").Value = 1961
$ws.Range("F14").Value = 935
$ws.Range("G14").Value = 4792
$ws.Range("H14").Value = 4423
$ws.Range("I14").Value = 18835
$ws.Range("J14").Value = 12075
$ws.Range("K14").Value = -9454
$ws.Range("L14").Value = 1168
$ws.Range("M14").Value = 3517
$ws.Range("N14").Value = 1112

$ws.Range("E15").Value = 55
$ws.Range("F15").Value = 164
$ws.Range("G15").Value = 308
$ws.Range("H15").Value = 398
$ws.Range("I15").Value = 246
$ws.Range("J15").Value = 434
$ws.Range("K15").Value = 284
$ws.Range("L15").Value = 323
$ws.Range("M15").Value = 317
$ws.Range("N15").Value = 772

$ws.Range("E16").Value = 3292
$ws.Range("F16").Value = 3307
$ws.Range("G16").Value = 3427
$ws.Range("H16").Value = 3514
$ws.Range("I16").Value = 3801
$ws.Range("J16").Value = 3946
$ws.Range("K16").Value = 4074
$ws.Range("L16").Value = 3992
$ws.Range("M16").Value = 4077
$ws.Range("N16").Value = 4083

$ws.Range("E17").Value = 120861
$ws.Range("F17").Value = 113176
$ws.Range("G17").Value = 167660
$ws.Range("H17").Value = 157247
$ws.Range("I17").Value = 188678
$ws.Range("J17").Value = 178440
$ws.Range("K17").Value = 287701
$ws.Range("L17").Value = 231598
$ws.Range("M17").Value = 252700
$ws.Range("N17").Value = 263443

$ws.Range("E19").Value = -19691
$ws.Range("F19").Value = 60204
$ws.Range("G19").Value = 32320
$ws.Range("H19").Value = 56227
$ws.Range("I19").Value = 24613
$ws.Range("J19").Value = 60190
$ws.Range("K19").Value = 75961
$ws.Range("L19").Value = 60432
$ws.Range("M19").Value = 78790
$ws.Range("N19").Value = 152147

$ws.Range("E20").Value = 248680
$ws.Range("F20").Value = 284219
$ws.Range("G20").Value = 285936
$ws.Range("H20").Value = 575366
$ws.Range("I20").Value = 345916
$ws.Range("J20").Value = 337659
$ws.Range("K20").Value = 384879
$ws.Range("L20").Value = 362269
$ws.Range("M20").Value = 493071
$ws.Range("N20").Value = 498764

$ws.Range("E26").Value = 261
$ws.Range("F26").Value = 292
$ws.Range("G26").Value = 266
$ws.Range("H26").Value = 294
$ws.Range("I26").Value = 287
$ws.Range("J26").Value = 297
$ws.Range("K26").Value = 288
$ws.Range("L26").Value = 282
$ws.Range("M26").Value = 272
$ws.Range("N26").Value = 276

$ws.Range("E27").Value = 718
$ws.Range("F27").Value = 716
$ws.Range("G27").Value = 716
$ws.Range("H27").Value = 725
$ws.Range("I27").Value = 734
$ws.Range("J27").Value = 743
$ws.Range("K27").Value = 733
$ws.Range("L27").Value = 724
$ws.Range("M27").Value = 717
$ws.Range("N27").Value = 711

